$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (rows 2-151).
# All of these cells currently store 45177 (2023-09-08) and should be bumped to
# 45178 (2023-09-09).
$ws.Range("C2:C151").Value = 45178
